$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (F) values for several rows
$ws.Range("F2").Value = 1
$ws.Range("F4").Value = -1
$ws.Range("F6").Value = -2
$ws.Range("F7").Value = -4
$ws.Range("F10").Value = 1
